$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-18 07:15:16"
$wsZh.Range("H4").Value = "2016-03-18 07:15:35"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-18 07:15:18"
$wsDe.Range("H4").Value = "2016-03-18 07:15:39"
